$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 1980.4615
$ws.Cells.Item(92, 9).Value = 1835.1
$ws.Cells.Item(92, 10).Value = 2465
$ws.Cells.Item(92, 11).Value = 1835.1
$ws.Cells.Item(92, 12).Value = 2465
$ws.Cells.Item(92, 13).Value = -587.0999999999999
$ws.Cells.Item(92, 14).Value = -4961
$ws.Cells.Item(103, 8).Value = 10085.444
$ws.Cells.Item(103, 9).Value = 731.2857
$ws.Cells.Item(103, 11).Value = 2193.8571
$ws.Cells.Item(103, 13).Value = -1607.8571
$ws.Cells.Item(129, 8).Value = 858.35
$ws.Cells.Item(129, 10).Value = 905.2198
$ws.Cells.Item(129, 12).Value = 2715.6594
$ws.Cells.Item(129, 14).Value = -12715.6594
$ws.Cells.Item(132, 8).Value = 307425.28
$ws.Cells.Item(132, 9).Value = 4749.3706
$ws.Cells.Item(132, 10).Value = 1669466.9
$ws.Cells.Item(132, 11).Value = 14248.1118
$ws.Cells.Item(132, 12).Value = 5008400.699999999
$ws.Cells.Item(132, 13).Value = -11718.1118
$ws.Cells.Item(132, 14).Value = -5013460.699999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(28, 8).Value = 3471
$ws.Cells.Item(28, 9).Value = 3471
$ws.Cells.Item(28, 11).Value = 3471
$ws.Cells.Item(28, 13).Value = -3279
$ws.Cells.Item(32, 8).Value = 5366.6064
$ws.Cells.Item(32, 9).Value = 4353.354
$ws.Cells.Item(32, 10).Value = 9107.846
$ws.Cells.Item(32, 11).Value = 4353.354
$ws.Cells.Item(32, 12).Value = 9107.846
$ws.Cells.Item(32, 13).Value = -4066.354
$ws.Cells.Item(32, 14).Value = -9681.846
$ws.Cells.Item(61, 8).Value = 2076
$ws.Cells.Item(61, 9).Value = 2052.75
$ws.Cells.Item(61, 10).Value = 2200
$ws.Cells.Item(61, 11).Value = 2052.75
$ws.Cells.Item(61, 12).Value = 2200
$ws.Cells.Item(61, 13).Value = -1840.75
$ws.Cells.Item(61, 14).Value = -2624
$ws.Cells.Item(74, 8).Value = 4673.1665
$ws.Cells.Item(74, 9).Value = 4647.8423
$ws.Cells.Item(74, 10).Value = 4769.4
$ws.Cells.Item(74, 11).Value = 4647.8423
$ws.Cells.Item(74, 12).Value = 4769.4
$ws.Cells.Item(74, 13).Value = -3773.8423
$ws.Cells.Item(74, 14).Value = -6517.4
$ws.Cells.Item(77, 8).Value = 4673.1665
$ws.Cells.Item(77, 9).Value = 4647.8423
$ws.Cells.Item(77, 10).Value = 4769.4
$ws.Cells.Item(77, 11).Value = 23239.2115
$ws.Cells.Item(77, 12).Value = 23847
$ws.Cells.Item(77, 13).Value = -18871.2115
$ws.Cells.Item(77, 14).Value = -32583
$ws.Cells.Item(94, 8).Value = 29132
$ws.Cells.Item(94, 10).Value = 29132
$ws.Cells.Item(94, 12).Value = 29132
$ws.Cells.Item(94, 14).Value = -30934
$ws.Cells.Item(99, 8).Value = 3471
$ws.Cells.Item(99, 9).Value = 3471
$ws.Cells.Item(99, 11).Value = 3471
$ws.Cells.Item(99, 13).Value = -476
$ws.Cells.Item(132, 8).Value = 1709.579
$ws.Cells.Item(132, 9).Value = 911
$ws.Cells.Item(132, 10).Value = 5968.6665
$ws.Cells.Item(132, 11).Value = 2733
$ws.Cells.Item(132, 12).Value = 17905.9995
$ws.Cells.Item(132, 13).Value = -203
$ws.Cells.Item(132, 14).Value = -22965.9995
$ws.Cells.Item(136, 8).Value = 2076
$ws.Cells.Item(136, 9).Value = 2052.75
$ws.Cells.Item(136, 10).Value = 2200
$ws.Cells.Item(136, 11).Value = 6158.25
$ws.Cells.Item(136, 12).Value = 6600
$ws.Cells.Item(136, 13).Value = -3608.25
$ws.Cells.Item(136, 14).Value = -11700

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2065.7334
$ws.Cells.Item(86, 9).Value = 1838.6842
$ws.Cells.Item(86, 10).Value = 2457.9092
$ws.Cells.Item(86, 11).Value = 1838.6842
$ws.Cells.Item(86, 12).Value = 2457.9092
$ws.Cells.Item(86, 13).Value = -715.6841999999999
$ws.Cells.Item(86, 14).Value = -4703.9092
$ws.Cells.Item(89, 8).Value = 2065.7334
$ws.Cells.Item(89, 9).Value = 1838.6842
$ws.Cells.Item(89, 10).Value = 2457.9092
$ws.Cells.Item(89, 11).Value = 9193.421
$ws.Cells.Item(89, 12).Value = 12289.546
$ws.Cells.Item(89, 13).Value = -3577.421
$ws.Cells.Item(89, 14).Value = -23521.546
$ws.Cells.Item(99, 8).Value = 5190.9
$ws.Cells.Item(99, 9).Value = 1710
$ws.Cells.Item(99, 11).Value = 1710
$ws.Cells.Item(99, 13).Value = -212
$ws.Cells.Item(134, 8).Value = 1902.119
$ws.Cells.Item(134, 9).Value = 1407.4103
$ws.Cells.Item(134, 10).Value = 8333.333000000001
$ws.Cells.Item(134, 11).Value = 4222.2309
$ws.Cells.Item(134, 12).Value = 24999.999
$ws.Cells.Item(134, 13).Value = -1687.2309
$ws.Cells.Item(134, 14).Value = -30069.999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 6174045
$ws.Cells.Item(16, 9).Value = 9260276
$ws.Cells.Item(16, 11).Value = 9260276
$ws.Cells.Item(16, 13).Value = -9259989
$ws.Cells.Item(107, 8).Value = 452.16666
$ws.Cells.Item(107, 9).Value = 395.03705
$ws.Cells.Item(107, 11).Value = 395.03705
$ws.Cells.Item(107, 13).Value = 1524.96295
$ws.Cells.Item(113, 8).Value = 6174045
$ws.Cells.Item(113, 9).Value = 9260276
$ws.Cells.Item(113, 11).Value = 9260276
$ws.Cells.Item(113, 13).Value = -9258106

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 246.66667
$ws.Cells.Item(33, 9).Value = 182.5
$ws.Cells.Item(33, 10).Value = 375
$ws.Cells.Item(33, 11).Value = 1095
$ws.Cells.Item(33, 12).Value = 2250
$ws.Cells.Item(33, 13).Value = -812
$ws.Cells.Item(33, 14).Value = -2816
$ws.Cells.Item(39, 8).Value = 12480.8
$ws.Cells.Item(39, 10).Value = 12480.8
$ws.Cells.Item(39, 12).Value = 37442.39999999999
$ws.Cells.Item(39, 14).Value = -38030.39999999999
$ws.Cells.Item(109, 8).Value = 3839.0908
$ws.Cells.Item(109, 9).Value = 3466.2
$ws.Cells.Item(109, 10).Value = 4149.8335
$ws.Cells.Item(109, 11).Value = 10398.6
$ws.Cells.Item(109, 12).Value = 12449.5005
$ws.Cells.Item(109, 13).Value = -9358.599999999999
$ws.Cells.Item(109, 14).Value = -14529.5005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6233.5
$ws.Cells.Item(70, 10).Value = 8207
$ws.Cells.Item(70, 12).Value = 8207
$ws.Cells.Item(70, 14).Value = -8747
$ws.Cells.Item(73, 8).Value = 6233.5
$ws.Cells.Item(73, 10).Value = 8207
$ws.Cells.Item(73, 12).Value = 8207
$ws.Cells.Item(73, 14).Value = -10079
$ws.Cells.Item(97, 8).Value = 928.5
$ws.Cells.Item(97, 9).Value = 925.6
$ws.Cells.Item(97, 10).Value = 933.3333
$ws.Cells.Item(97, 11).Value = 925.6
$ws.Cells.Item(97, 12).Value = 933.3333
$ws.Cells.Item(97, 13).Value = -429.6
$ws.Cells.Item(97, 14).Value = -1925.3333
$ws.Cells.Item(102, 8).Value = 2957.5386
$ws.Cells.Item(102, 9).Value = 2033.8
$ws.Cells.Item(102, 10).Value = 6036.6665
$ws.Cells.Item(102, 11).Value = 2033.8
$ws.Cells.Item(102, 12).Value = 6036.6665
$ws.Cells.Item(102, 13).Value = -411.8
$ws.Cells.Item(102, 14).Value = -9280.666499999999
$ws.Cells.Item(132, 8).Value = 2787.182
$ws.Cells.Item(132, 9).Value = 1704.421
$ws.Cells.Item(132, 10).Value = 4256.643
$ws.Cells.Item(132, 11).Value = 5113.263
$ws.Cells.Item(132, 12).Value = 12769.929
$ws.Cells.Item(132, 13).Value = -2583.263
$ws.Cells.Item(132, 14).Value = -17829.929

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3030.182
$ws.Cells.Item(46, 9).Value = 3356
$ws.Cells.Item(46, 10).Value = 2758.6667
$ws.Cells.Item(46, 11).Value = 3356
$ws.Cells.Item(46, 12).Value = 2758.6667
$ws.Cells.Item(46, 13).Value = -3168
$ws.Cells.Item(46, 14).Value = -3134.6667
$ws.Cells.Item(74, 8).Value = 43669.5
$ws.Cells.Item(74, 10).Value = 43669.5
$ws.Cells.Item(74, 12).Value = 43669.5
$ws.Cells.Item(74, 14).Value = -45665.5
$ws.Cells.Item(77, 8).Value = 43669.5
$ws.Cells.Item(77, 10).Value = 43669.5
$ws.Cells.Item(77, 12).Value = 131008.5
$ws.Cells.Item(77, 14).Value = -140992.5
$ws.Cells.Item(92, 8).Value = 32759.334
$ws.Cells.Item(92, 10).Value = 32759.334
$ws.Cells.Item(92, 12).Value = 32759.334
$ws.Cells.Item(92, 14).Value = -37751.334
$ws.Cells.Item(136, 8).Value = 3855.0833
$ws.Cells.Item(136, 9).Value = 1495.125
$ws.Cells.Item(136, 10).Value = 8575
$ws.Cells.Item(136, 11).Value = 4485.375
$ws.Cells.Item(136, 12).Value = 25725
$ws.Cells.Item(136, 13).Value = -1935.375
$ws.Cells.Item(136, 14).Value = -30825

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2000
$ws.Cells.Item(81, 9).Value = 2000
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 4000
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = -2939
$ws.Cells.Item(81, 14).ClearContents()
$ws.Cells.Item(84, 8).Value = 2000
$ws.Cells.Item(84, 9).Value = 2000
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 20000
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = -14696
$ws.Cells.Item(84, 14).ClearContents()
$ws.Cells.Item(94, 8).Value = 10330
$ws.Cells.Item(94, 10).Value = 10330
$ws.Cells.Item(94, 12).Value = 10330
$ws.Cells.Item(94, 14).Value = -12132
